$d = $word.ActiveDocument

# 1. Fix typo: "usual categories or video games" -> "usual categories of video games"
$rng = $d.Content
$found = $rng.Find.Execute("usual categories or video games", $true, $false, $false, $false, $false, $true, 1, $false, "usual categories of video games", 2)
if (-not $found) { throw "Step 1 (typo fix) text not found" }

# 2. "mobile-responsive (text is too small)" -> "mobile-responsive (the text is too small)"
$rng = $d.Content
$found = $rng.Find.Execute("mobile-responsive (text is too small)", $true, $false, $false, $false, $false, $true, 1, $false, "mobile-responsive (the text is too small)", 2)
if (-not $found) { throw "Step 2 (mobile-responsive) text not found" }

# 3. "somewhat close the color all Octonauts wear" -> "somewhat close to the color all Octonauts wear"
$rng = $d.Content
$found = $rng.Find.Execute("somewhat close the color all Octonauts wear", $true, $false, $false, $false, $false, $true, 1, $false, "somewhat close to the color all Octonauts wear", 2)
if (-not $found) { throw "Step 3 (close to) text not found" }

# 4. "Safari, Chrome, Firefox for viewing and testing." -> "Safari (OSX and iOS), Chrome, Firefox for viewing and testing."
$rng = $d.Content
$found = $rng.Find.Execute("Safari, Chrome, Firefox for viewing and testing.", $true, $false, $false, $false, $false, $true, 1, $false, "Safari (OSX and iOS), Chrome, Firefox for viewing and testing.", 2)
if (-not $found) { throw "Step 4 (Safari) text not found" }

# 5. "Bootstrap stylesheet." -> "Bootstrap stylesheet for some CSS."
$rng = $d.Content
$found = $rng.Find.Execute("Bootstrap stylesheet.", $true, $false, $false, $false, $false, $true, 1, $false, "Bootstrap stylesheet for some CSS.", 2)
if (-not $found) { throw "Step 5 (Bootstrap stylesheet) text not found" }

# 6. "were probably my most used resource" -> "were by far my most used resource"
#    (the word "probably" is swapped out for "by far")
$rng = $d.Content
$found = $rng.Find.Execute("were probably my most used resource", $true, $false, $false, $false, $false, $true, 1, $false, "were by far my most used resource", 2)
if (-not $found) { throw "Step 6 (by far) text not found" }

# Reposition the _GoBack bookmark (Word always tracks the most-recent-edit location
# with this hidden bookmark) so that it sits right after the newly inserted "by far".
$bfRng = $d.Content
$found = $bfRng.Find.Execute("by far", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not locate 'by far' to reposition _GoBack bookmark" }

$gobackPoint = $d.Range($bfRng.End, $bfRng.End)
$d.Bookmarks.Add("_GoBack", $gobackPoint)
